# VT_ELAND_ALL_V01.xlsx -- "Add files via upload" commit
#
# 1) Rename sheet "Arkusz1" -> "GRID" (formulas referencing it update automatically).
# 2) Update the small FI_T table on the GRID sheet (headers/units row + new "-" cell,
#    extend formatting down one row).
# 3) Re-point each sheet's remembered cell selection to match the new file, and make
#    GRID the active/selected tab (it was Demand before).

$wb = $excel.ActiveWorkbook

# --- 1. Rename Arkusz1 -> GRID -------------------------------------------------
$grid = $wb.Worksheets.Item("Arkusz1")
$grid.Name = "GRID"

# --- 2. Fix up the GRID table contents -----------------------------------------
# Row 6 header swap: G6 CAP2ACT, H6 Efficiency (was G6 Efficiency, H6 Extraction cost)
$grid.Range("G6").Value = "CAP2ACT"
$grid.Range("H6").Value = "Efficiency"

# Row 7 units swap: F7 cleared, H7 becomes a quote-prefixed "-" placeholder
$grid.Range("F7").ClearContents()
$grid.Range("H7").Value = "'-"

# Extend the row-7 side formatting down into new row 8/9 border cells
$grid.Range("B7").Copy()
$grid.Range("B8:B9").PasteSpecial(-4122)
$grid.Range("I7").Copy()
$grid.Range("I8:I9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Selections / active tab -------------------------------------------------
$wb.Worksheets.Item("FI_Comm").Range("F14").Select()
$wb.Worksheets.Item("FI_Process").Range("C17").Select()
$wb.Worksheets.Item("Power Plants").Range("I23").Select()
$wb.Worksheets.Item("Demand").Range("J7").Select()

# GRID becomes the active tab last, matching the saved workbook view state.
$grid.Range("I12").Select()
